$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 8.240429235646172
$ws.Range("C2").Value = 4.443077152129236
$ws.Range("D2").Value = 11.25382604485561
$ws.Range("F2").Value = 30.23442155288505
$ws.Range("G2").Value = 28.57673330631639
$ws.Range("H2").Value = 14.35225666546571
$ws.Range("I2").Value = 20.61435096332924
$ws.Range("J2").Value = 11.41980662888941
$ws.Range("K2").Value = 8.192135282892176
$ws.Range("M2").Value = 15.62145121917271
$ws.Range("N2").Value = 18.64775362942765
$ws.Range("O2").Value = 21.80188204576866
# Row 3
$ws.Range("B3").Value = 7.956720529763678
$ws.Range("C3").Value = 4.276279489553058
$ws.Range("D3").Value = 11.21024531140881
$ws.Range("F3").Value = 30.28823827716163
$ws.Range("G3").Value = 28.65007657822885
$ws.Range("H3").Value = 14.39364635872947
$ws.Range("I3").Value = 20.6930061555147
$ws.Range("J3").Value = 11.43783597741448
$ws.Range("K3").Value = 8.007157045993043
$ws.Range("M3").Value = 15.54049077128654
$ws.Range("N3").Value = 18.69849863201008
$ws.Range("O3").Value = 21.87004078790364
# Row 4
$ws.Range("B4").Value = 7.778209771711259
$ws.Range("C4").Value = 4.169694096498313
$ws.Range("D4").Value = 11.18545720380998
$ws.Range("F4").Value = 30.3276890096703
$ws.Range("G4").Value = 28.70281298828469
$ws.Range("H4").Value = 14.42095960690261
$ws.Range("I4").Value = 20.74461671573278
$ws.Range("J4").Value = 11.45053964519616
$ws.Range("K4").Value = 7.892320114048172
$ws.Range("M4").Value = 15.49302434825981
$ws.Range("N4").Value = 18.73118577248602
$ws.Range("O4").Value = 21.91578071670187
# Row 5
$ws.Range("B5").Value = 7.704492357893569
$ws.Range("C5").Value = 4.125251848292014
$ws.Range("D5").Value = 11.17585925672032
$ws.Range("F5").Value = 30.34537500756615
$ws.Range("G5").Value = 28.72623512839976
$ws.Range("H5").Value = 14.43256806582291
$ws.Range("I5").Value = 20.76648258540046
$ws.Range("J5").Value = 11.45612743068343
$ws.Range("K5").Value = 7.84527063299338
$ws.Range("M5").Value = 15.47426068421694
$ws.Range("N5").Value = 18.74489168134258
$ws.Range("O5").Value = 21.9353976302963
# Row 6
$ws.Range("B6").Value = 7.692196341338507
$ws.Range("C6").Value = 4.11781267360401
$ws.Range("D6").Value = 11.17429614844693
$ws.Range("F6").Value = 30.34840890959043
$ws.Range("G6").Value = 28.73024085294326
$ws.Range("H6").Value = 14.43452452778058
$ws.Range("I6").Value = 20.77016378676593
$ws.Range("J6").Value = 11.45708010165514
$ws.Range("K6").Value = 7.837444778954092
$ws.Range("M6").Value = 15.47118040652201
$ws.Range("N6").Value = 18.74719085868884
$ws.Range("O6").Value = 21.9387140227398
# Row 7
$ws.Range("B7").Value = 7.777219376926463
$ws.Range("C7").Value = 4.169098755612263
$ws.Range("D7").Value = 11.18532571427519
$ws.Range("F7").Value = 30.32792101480896
$ws.Range("G7").Value = 28.70312105396462
$ws.Range("H7").Value = 14.42111422638851
$ws.Range("I7").Value = 20.74490822862026
$ws.Range("J7").Value = 11.45061333997191
$ws.Range("K7").Value = 7.891686523914726
$ws.Range("M7").Value = 15.49276893000108
$ws.Range("N7").Value = 18.73136905238265
$ws.Range("O7").Value = 21.91604132025731
# Row 8
$ws.Range("B8").Value = 8.143565317428846
$ws.Range("C8").Value = 4.386454045881691
$ws.Range("D8").Value = 11.23839529482203
$ws.Range("F8").Value = 30.25164692901718
$ws.Range("G8").Value = 28.60042011930922
$ws.Range("H8").Value = 14.36613374706932
$ws.Range("I8").Value = 20.64078324095866
$ws.Range("J8").Value = 11.42568420981998
$ws.Range("K8").Value = 8.128653917482522
$ws.Range("M8").Value = 15.5930796814471
$ws.Range("N8").Value = 18.66493367838433
$ws.Range("O8").Value = 21.82457536372381
# Row 9
$ws.Range("B9").Value = 8.823363460338715
$ws.Range("C9").Value = 4.778038726579395
$ws.Range("D9").Value = 11.35772994621868
$ws.Range("F9").Value = 30.15296790180046
$ws.Range("G9").Value = 28.46036788667466
$ws.Range("H9").Value = 14.27337488055309
$ws.Range("I9").Value = 20.46289031019127
$ws.Range("J9").Value = 11.38975175352735
$ws.Range("K9").Value = 8.580671588847473
$ws.Range("M9").Value = 15.80689464974914
$ws.Range("N9").Value = 18.54674252332213
$ws.Range("O9").Value = 21.67610798335565
# Row 10
$ws.Range("B10").Value = 9.294098665298741
$ws.Range("C10").Value = 5.042846402687061
$ws.Range("D10").Value = 11.45418039562529
$ws.Range("F10").Value = 30.11155149501011
$ws.Range("G10").Value = 28.39513772960638
$ws.Range("H10").Value = 14.21438184666487
$ws.Range("I10").Value = 20.34820242390444
$ws.Range("J10").Value = 11.37123628945654
$ws.Range("K10").Value = 8.90167300436498
$ws.Range("M10").Value = 15.97343729890978
$ws.Range("N10").Value = 18.46720950652885
$ws.Range("O10").Value = 21.58590444763876
# Row 11
$ws.Range("B11").Value = 9.501107831471058
$ws.Range("C11").Value = 5.158049892435758
$ws.Range("D11").Value = 11.49984042257962
$ws.Range("F11").Value = 30.09946257875157
$ws.Range("G11").Value = 28.37368674236592
$ws.Range("H11").Value = 14.18952800351153
$ws.Range("I11").Value = 20.29950044973549
$ws.Range("J11").Value = 11.36452137666526
$ws.Range("K11").Value = 9.044650919379889
$ws.Range("M11").Value = 16.05103487937258
$ws.Range("N11").Value = 18.43259912971821
$ws.Range("O11").Value = 21.5489749824316
# Row 12
$ws.Range("B12").Value = 9.578407549519765
$ws.Range("C12").Value = 5.200898158938204
$ws.Range("D12").Value = 11.51737617303991
$ws.Range("F12").Value = 30.09585507637684
$ws.Range("G12").Value = 28.36674877941378
$ws.Range("H12").Value = 14.18040123023472
$ws.Range("I12").Value = 20.28155710274273
$ws.Range("J12").Value = 11.36222378311079
$ws.Range("K12").Value = 9.098305892132965
$ws.Range("M12").Value = 16.08066391023345
$ws.Range("N12").Value = 18.4197177369782
$ws.Range("O12").Value = 21.53558159991904
# Row 13
$ws.Range("B13").Value = 9.56180910532113
$ws.Range("C13").Value = 5.191704862303745
$ws.Range("D13").Value = 11.51358880980253
$ws.Range("F13").Value = 30.09658887621377
$ws.Range("G13").Value = 28.36819025301505
$ws.Range("H13").Value = 14.18235417950052
$ws.Range("I13").Value = 20.28539932811652
$ws.Range("J13").Value = 11.36270771260591
$ws.Range("K13").Value = 9.086772771057746
$ws.Range("M13").Value = 16.0742722091073
$ws.Range("N13").Value = 18.42248199128299
$ws.Range("O13").Value = 21.53843980766316
# Row 14
$ws.Range("B14").Value = 9.507489517301602
$ws.Range("C14").Value = 5.161590720500412
$ws.Range("D14").Value = 11.50127824933961
$ws.Range("F14").Value = 30.09914634596966
$ws.Range("G14").Value = 28.37309218882309
$ws.Range("H14").Value = 14.1887714310955
$ws.Range("I14").Value = 20.29801423812154
$ws.Range("J14").Value = 11.36432744057721
$ws.Range("K14").Value = 9.049075189976463
$ws.Range("M14").Value = 16.05346769061225
$ws.Range("N14").Value = 18.43153486945206
$ws.Range("O14").Value = 21.54786125194453
# Row 15
$ws.Range("B15").Value = 9.474073442367557
$ws.Range("C15").Value = 5.143043217909591
$ws.Range("D15").Value = 11.49376927811575
$ws.Range("F15").Value = 30.10083920741205
$ws.Range("G15").Value = 28.37624916307288
$ws.Range("H15").Value = 14.19273927059962
$ws.Range("I15").Value = 20.30580622507928
$ws.Range("J15").Value = 11.36535149013968
$ws.Range("K15").Value = 9.025919420007293
$ws.Range("M15").Value = 16.04075557770666
$ws.Range("N15").Value = 18.43710926858369
$ws.Range("O15").Value = 21.5537091468581
# Row 16
$ws.Range("B16").Value = 9.280420922567686
$ws.Range("C16").Value = 5.035209933699398
$ws.Range("D16").Value = 11.45123135338648
$ws.Range("F16").Value = 30.11247735688876
$ws.Range("G16").Value = 28.39670530478991
$ws.Range("H16").Value = 14.21604597494712
$ws.Range("I16").Value = 20.35145504257583
$ws.Range("J16").Value = 11.37170945939079
$ws.Range("K16").Value = 8.892263429903794
$ws.Range("M16").Value = 15.96840140484524
$ws.Range("N16").Value = 18.46950288308021
$ws.Range("O16").Value = 21.58840053670695
# Row 17
$ws.Range("B17").Value = 9.159746395881424
$ws.Range("C17").Value = 4.967695171236269
$ws.Range("D17").Value = 11.42558479851745
$ws.Range("F17").Value = 30.12134599071
$ws.Range("G17").Value = 28.41136259599798
$ws.Range("H17").Value = 14.23085145690229
$ws.Range("I17").Value = 20.38034791190958
$ws.Range("J17").Value = 11.37604703392359
$ws.Range("K17").Value = 8.809452967512096
$ws.Range("M17").Value = 15.92447102932667
$ws.Range("N17").Value = 18.48977668261796
$ws.Range("O17").Value = 21.6107344259716
# Row 18
$ws.Range("B18").Value = 9.089669772162006
$ws.Range("C18").Value = 4.928368502280925
$ws.Range("D18").Value = 11.41100222037906
$ws.Range("F18").Value = 30.12708255356282
$ws.Range("G18").Value = 28.42056691067969
$ws.Range("H18").Value = 14.23955378321002
$ws.Range("I18").Value = 20.39729294679447
$ws.Range("J18").Value = 11.37870267155335
$ws.Range("K18").Value = 8.761538482975954
$ws.Range("M18").Value = 15.89937749663667
$ws.Range("N18").Value = 18.50158542629176
$ws.Range("O18").Value = 21.62396658743053
# Row 19
$ws.Range("B19").Value = 9.065830436760322
$ws.Range("C19").Value = 4.914969006323513
$ws.Range("D19").Value = 11.40609411167891
$ws.Range("F19").Value = 30.12913402217135
$ws.Range("G19").Value = 28.42381614645379
$ws.Range("H19").Value = 14.24253229786926
$ws.Range("I19").Value = 20.40308634119062
$ws.Range("J19").Value = 11.37962944976891
$ws.Range("K19").Value = 8.745268235768815
$ws.Range("M19").Value = 15.89091174231107
$ws.Range("N19").Value = 18.50560907474667
$ws.Range("O19").Value = 21.62851309284802
# Row 20
$ws.Range("B20").Value = 9.172662017410941
$ws.Range("C20").Value = 4.974933514940324
$ws.Range("D20").Value = 11.4282975447493
$ws.Range("F20").Value = 30.12033613401686
$ws.Range("G20").Value = 28.40972219241665
$ws.Range("H20").Value = 14.2292560768397
$ws.Range("I20").Value = 20.37723841245682
$ws.Range("J20").Value = 11.37556865427576
$ws.Range("K20").Value = 8.818298055620881
$ws.Range("M20").Value = 15.92912961911042
$ws.Range("N20").Value = 18.48760321284551
$ws.Range("O20").Value = 21.6083169590767
# Row 21
$ws.Range("B21").Value = 9.523474560895336
$ws.Range("C21").Value = 5.170457205119805
$ws.Range("D21").Value = 11.50488759408394
$ws.Range("F21").Value = 30.09836882822987
$ws.Range("G21").Value = 28.37162019148171
$ws.Range("H21").Value = 14.18687880012207
$ws.Range("I21").Value = 20.29429538924374
$ws.Range("J21").Value = 11.36384503601581
$ws.Range("K21").Value = 9.060161494279564
$ws.Range("M21").Value = 16.0595720038335
$ws.Range("N21").Value = 18.42886972563433
$ws.Range("O21").Value = 21.54507789979421
# Row 22
$ws.Range("B22").Value = 9.746371922974632
$ws.Range("C22").Value = 5.293707977766809
$ws.Range("D22").Value = 11.55636785873034
$ws.Range("F22").Value = 30.08966700231532
$ws.Range("G22").Value = 28.35362650551414
$ws.Range("H22").Value = 14.16084294932561
$ws.Range("I22").Value = 20.2429961111406
$ws.Range("J22").Value = 11.35761195205709
$ws.Range("K22").Value = 9.215370422562179
$ws.Range("M22").Value = 16.14623952245128
$ws.Range("N22").Value = 18.39179392769704
$ws.Range("O22").Value = 21.50719266631818
# Row 23
$ws.Range("B23").Value = 9.628010472764036
$ws.Range("C23").Value = 5.228347739553114
$ws.Range("D23").Value = 11.52876537820185
$ws.Range("F23").Value = 30.09379420752491
$ws.Range("G23").Value = 28.36259727750112
$ws.Range("H23").Value = 14.17458694319213
$ws.Range("I23").Value = 20.27010933016689
$ws.Range("J23").Value = 11.36080805948461
$ws.Range("K23").Value = 9.132809835806061
$ws.Range("M23").Value = 16.09986044022141
$ws.Range("N23").Value = 18.41146241419783
$ws.Range("O23").Value = 21.52709725572166
# Row 24
$ws.Range("B24").Value = 9.166825037861523
$ws.Range("C24").Value = 4.971662651487511
$ws.Range("D24").Value = 11.42707060681946
$ws.Range("F24").Value = 30.12079070329486
$ws.Range("G24").Value = 28.41046139679982
$ws.Range("H24").Value = 14.22997675488632
$ws.Range("I24").Value = 20.3786431763432
$ws.Range("J24").Value = 11.3757844253601
$ws.Range("K24").Value = 8.814300139505429
$ws.Range("M24").Value = 15.92702296262587
$ws.Range("N24").Value = 18.48858536166766
$ws.Range("O24").Value = 21.60940867452303
# Row 25
$ws.Range("B25").Value = 8.644164993768758
$ws.Range("C25").Value = 4.67601245497753
$ws.Range("D25").Value = 11.32386528468359
$ws.Range("F25").Value = 30.17420698717678
$ws.Range("G25").Value = 28.49165821192757
$ws.Range("H25").Value = 14.29685894232266
$ws.Range("I25").Value = 20.50820213070641
$ws.Range("J25").Value = 11.39808676888295
$ws.Range("K25").Value = 8.580671588847473
$ws.Range("M25").Value = 15.74731901018369
$ws.Range("N25").Value = 18.57742922458306
$ws.Range("O25").Value = 21.71296033555238
